# Applies the "Updated cryptos list" data refresh described in the commit diff.
# Updates Price (col D) and Volume(1h) (col E) values for the changed rows,
# and, for three rows, the Coin name / Link that got reordered in the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '51.972.90'
$ws.Cells.Item(2, 5).Value = '  -0.78%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.837.37'
$ws.Cells.Item(3, 5).Value = '  +0.65%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.03%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '354.17'
$ws.Cells.Item(5, 5).Value = '  +4.94%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '113.08'
$ws.Cells.Item(6, 5).Value = '  -4.31%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.561'
$ws.Cells.Item(7, 5).Value = '  +3.74%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.604'
$ws.Cells.Item(9, 5).Value = '  +3.86%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '41.47'
$ws.Cells.Item(10, 5).Value = '  -2.60%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0855'
$ws.Cells.Item(11, 5).Value = '  -1.40%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  +0.83%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '19.90'
$ws.Cells.Item(13, 5).Value = '  -2.86%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '7.74'
$ws.Cells.Item(14, 5).Value = '  +0.31%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '3.276.96'
$ws.Cells.Item(15, 5).Value = '  +0.54%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '2.831.51'
$ws.Cells.Item(16, 5).Value = '  -0.16%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.891'
$ws.Cells.Item(17, 5).Value = '  -0.37%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '51.873.31'
$ws.Cells.Item(18, 5).Value = '  -1.03%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '7.60'
$ws.Cells.Item(19, 5).Value = '  +8.86%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '3.21'
$ws.Cells.Item(20, 5).Value = '  -2.60%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '13.44'
$ws.Cells.Item(21, 5).Value = '  -1.80%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '0.0₃0985'
$ws.Cells.Item(22, 5).Value = '  +0.39%  '

# Row 23
$ws.Cells.Item(23, 2).Value = 'Litecoin'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '69.71'
$ws.Cells.Item(23, 5).Value = '  -1.15%  '

# Row 24
$ws.Cells.Item(24, 2).Value = 'BitcoinCash'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '267.81'
$ws.Cells.Item(24, 5).Value = '  -4.97%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  +1.84%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '26.70'
$ws.Cells.Item(26, 5).Value = '  -1.77%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '1.00'
$ws.Cells.Item(27, 5).Value = '  -0.05%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '10.38'
$ws.Cells.Item(28, 5).Value = '  +0.64%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  +0.76%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  -3.33%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '51.73'
$ws.Cells.Item(31, 5).Value = '  +2.35%  '

# Row 32
$ws.Cells.Item(32, 2).Value = 'VeChain'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.0461'
$ws.Cells.Item(32, 5).Value = '  +29.12%  '

# Row 33
$ws.Cells.Item(33, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '33.72'
$ws.Cells.Item(33, 5).Value = '  -6.37%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '5.85'
$ws.Cells.Item(34, 5).Value = '  +3.25%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '5.24'
$ws.Cells.Item(35, 5).Value = '  +3.75%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.0829'
$ws.Cells.Item(36, 5).Value = '  -0.19%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.999'
$ws.Cells.Item(37, 5).Value = '  +0.00%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  -4.36%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '3.22'
$ws.Cells.Item(39, 5).Value = '  -3.88%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  -6.03%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '23.51'
$ws.Cells.Item(41, 5).Value = '  -0.95%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.116'
$ws.Cells.Item(42, 5).Value = '  +1.36%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '2.53'
$ws.Cells.Item(43, 5).Value = '  -6.44%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '125.30'
$ws.Cells.Item(44, 5).Value = '  -2.38%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  -3.06%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'Maker'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(46, 4).Value = '2.103.56'
$ws.Cells.Item(46, 5).Value = '  -0.65%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'NEARProtocol'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '3.35'
$ws.Cells.Item(47, 5).Value = '  -0.57%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '2.26'
$ws.Cells.Item(48, 5).Value = '  +0.77%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '5.83'
$ws.Cells.Item(49, 5).Value = '  +4.05%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.959'
$ws.Cells.Item(50, 5).Value = '  +3.84%  '

# Row 51
$ws.Cells.Item(51, 2).Value = 'FraxShare'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '8.93'
$ws.Cells.Item(51, 5).Value = '  -0.22%  '
